$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '27.447.40'
$cell.Style = 'Normal'
$ws.Range('E2').Value = '  -2.93%  '

$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '1.744.72'
$cell.Style = 'Normal'
$ws.Range('E3').Value = '  -3.31%  '

$ws.Range('E4').Value = '  +0.14%  '

$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '322.62'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  -4.62%  '

$ws.Range('E6').Value = '  +0.09%  '

$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '0.4220'
$cell.Style = 'Normal'
$ws.Range('E7').Value = '  -9.80%  '

$ws.Range('E8').Value = '  -5.69%  '

$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '45.40'
$cell.Style = 'Normal'
$ws.Range('E9').Value = '  -0.25%  '

$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '0.07407'
$cell.Style = 'Normal'
$ws.Range('E10').Value = '  -2.66%  '

$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '1.110'
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  -3.44%  '

$ws.Range('E12').Value = '  +0.08%  '

$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '21.38'
$cell.Style = 'Normal'
$ws.Range('E13').Value = '  -4.40%  '

$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '6.087'
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  -3.75%  '

$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '7.173'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  -3.69%  '

$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '1.741.70'
$cell.Style = 'Normal'
$ws.Range('E16').Value = '  -3.50%  '

$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '0.00001063'
$cell.Style = 'Normal'
$ws.Range('E17').Value = '  -2.74%  '

$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '87.27'
$cell.Style = 'Normal'
$ws.Range('E18').Value = '  +6.77%  '

$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '0.06184'
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  -8.08%  '

$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '1.001'
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  +0.15%  '

$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '16.82'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  -3.40%  '

$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '6.089'
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  -4.91%  '

$ws.Range('E23').Value = '  -5.48%  '

$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '27.482.45'
$cell.Style = 'Normal'
$ws.Range('E24').Value = '  -2.73%  '

$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '11.59'
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  -2.09%  '

$ws.Range('E26').Value = '  -3.46%  '

$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '20.42'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  -1.42%  '

$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '151.80'
$cell.Style = 'Normal'
$ws.Range('E28').Value = '  -1.38%  '

$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '2.356'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  -0.54%  '

$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '1.939.62'
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  -3.50%  '

$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '126.03'
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  -5.28%  '

$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '1.204'
$cell.Style = 'Normal'
$ws.Range('E32').Value = '  -3.82%  '

$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '5.667'
$cell.Style = 'Normal'

$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '3.689'
$cell.Style = 'Normal'
$ws.Range('E34').Value = '  -8.63%  '

$ws.Range('B35').Value = 'Stellar'
$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '0.09126'
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  -5.18%  '

$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '12.60'
$cell.Style = 'Normal'
$ws.Range('E36').Value = '  +4.13%  '

$ws.Range('E37').Value = '  -2.90%  '

$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '0.2126'
$cell.Style = 'Normal'
$ws.Range('E38').Value = '  -4.88%  '

$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '5.085'
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  -3.13%  '

$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '0.06072'
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  -4.67%  '

$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '0.6386'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  -3.55%  '

$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '1.194'
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  -3.46%  '

$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '1.423'
$cell.Style = 'Normal'
$ws.Range('E43').Value = '  -5.25%  '

$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '1.001'
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  +0.08%  '

$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '7.875'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  -4.29%  '

$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '13.67'
$cell.Style = 'Normal'
$ws.Range('E46').Value = '  -3.59%  '

$ws.Range('E47').Value = '  -3.24%  '

$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '0.5853'
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  -4.35%  '

$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '124.79'
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  -4.38%  '

$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '1.947'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  -4.22%  '

$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '0.06852'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  -4.30%  '

